$d = $word.ActiveDocument

$d.Content.Find.Execute("50-4=46", $true, $false, $false, $false, $false, $true, 1, $false, "8+41=49", 2) | Out-Null
$d.Content.Find.Execute("5+77=82", $true, $false, $false, $false, $false, $true, 1, $false, "89-59=30", 2) | Out-Null
$d.Content.Find.Execute("71+11=82", $true, $false, $false, $false, $false, $true, 1, $false, "30+28=58", 2) | Out-Null
$d.Content.Find.Execute("65-54=11", $true, $false, $false, $false, $false, $true, 1, $false, "17+3=20", 2) | Out-Null
$d.Content.Find.Execute("85+5=90", $true, $false, $false, $false, $false, $true, 1, $false, "98-22=76", 2) | Out-Null
$d.Content.Find.Execute("22+8=30", $true, $false, $false, $false, $false, $true, 1, $false, "27+12=39", 2) | Out-Null
$d.Content.Find.Execute("62+14=76", $true, $false, $false, $false, $false, $true, 1, $false, "16+2=18", 2) | Out-Null
$d.Content.Find.Execute("33+26=59", $true, $false, $false, $false, $false, $true, 1, $false, "44-23=21", 2) | Out-Null
$d.Content.Find.Execute("92-15=77", $true, $false, $false, $false, $false, $true, 1, $false, "9+78=87", 2) | Out-Null
$d.Content.Find.Execute("36-19=17", $true, $false, $false, $false, $false, $true, 1, $false, "47+26=73", 2) | Out-Null
$d.Content.Find.Execute("80-34=46", $true, $false, $false, $false, $false, $true, 1, $false, "27+53=80", 2) | Out-Null
$d.Content.Find.Execute("72-41=31", $true, $false, $false, $false, $false, $true, 1, $false, "61-22=39", 2) | Out-Null
$d.Content.Find.Execute("94-69=25", $true, $false, $false, $false, $false, $true, 1, $false, "24+59=83", 2) | Out-Null
$d.Content.Find.Execute("69-50=19", $true, $false, $false, $false, $false, $true, 1, $false, "90-76=14", 2) | Out-Null
$d.Content.Find.Execute("50-35=15", $true, $false, $false, $false, $false, $true, 1, $false, "12+16=28", 2) | Out-Null
$d.Content.Find.Execute("70-46=24", $true, $false, $false, $false, $false, $true, 1, $false, "5+67=72", 2) | Out-Null
$d.Content.Find.Execute("42+10=52", $true, $false, $false, $false, $false, $true, 1, $false, "64+26=90", 2) | Out-Null
$d.Content.Find.Execute("10+59=69", $true, $false, $false, $false, $false, $true, 1, $false, "80-54=26", 2) | Out-Null
$d.Content.Find.Execute("9+79=88", $true, $false, $false, $false, $false, $true, 1, $false, "87-28=59", 2) | Out-Null
$d.Content.Find.Execute("3+13=16", $true, $false, $false, $false, $false, $true, 1, $false, "7+3=10", 2) | Out-Null
$d.Content.Find.Execute("71-28=43", $true, $false, $false, $false, $false, $true, 1, $false, "37+58=95", 2) | Out-Null
$d.Content.Find.Execute("53-17=36", $true, $false, $false, $false, $false, $true, 1, $false, "19+37=56", 2) | Out-Null
$d.Content.Find.Execute("96-32=64", $true, $false, $false, $false, $false, $true, 1, $false, "35+34=69", 2) | Out-Null
$d.Content.Find.Execute("14+31=45", $true, $false, $false, $false, $false, $true, 1, $false, "44-17=27", 2) | Out-Null
$d.Content.Find.Execute("37+26=63", $true, $false, $false, $false, $false, $true, 1, $false, "13-6=7", 2) | Out-Null
$d.Content.Find.Execute("29-14=15", $true, $false, $false, $false, $false, $true, 1, $false, "62-45=17", 2) | Out-Null
$d.Content.Find.Execute("50+8=58", $true, $false, $false, $false, $false, $true, 1, $false, "81+0=81", 2) | Out-Null
$d.Content.Find.Execute("15+44=59", $true, $false, $false, $false, $false, $true, 1, $false, "76-71=5", 2) | Out-Null
$d.Content.Find.Execute("48+25=73", $true, $false, $false, $false, $false, $true, 1, $false, "21+11=32", 2) | Out-Null
$d.Content.Find.Execute("43+44=87", $true, $false, $false, $false, $false, $true, 1, $false, "30+21=51", 2) | Out-Null
$d.Content.Find.Execute("14+84=98", $true, $false, $false, $false, $false, $true, 1, $false, "50+21=71", 2) | Out-Null
$d.Content.Find.Execute("18+68=86", $true, $false, $false, $false, $false, $true, 1, $false, "11+84=95", 2) | Out-Null
$d.Content.Find.Execute("44-32=12", $true, $false, $false, $false, $false, $true, 1, $false, "56+13=69", 2) | Out-Null
$d.Content.Find.Execute("56-4=52", $true, $false, $false, $false, $false, $true, 1, $false, "41-11=30", 2) | Out-Null
$d.Content.Find.Execute("68-43=25", $true, $false, $false, $false, $false, $true, 1, $false, "9+44=53", 2) | Out-Null
$d.Content.Find.Execute("94-89=5", $true, $false, $false, $false, $false, $true, 1, $false, "17+70=87", 2) | Out-Null
$d.Content.Find.Execute("24+49=73", $true, $false, $false, $false, $false, $true, 1, $false, "5+2=7", 2) | Out-Null
$d.Content.Find.Execute("73+10=83", $true, $false, $false, $false, $false, $true, 1, $false, "77+19=96", 2) | Out-Null
$d.Content.Find.Execute("40-21=19", $true, $false, $false, $false, $false, $true, 1, $false, "49+35=84", 2) | Out-Null
$d.Content.Find.Execute("90-50=40", $true, $false, $false, $false, $false, $true, 1, $false, "95-62=33", 2) | Out-Null
$d.Content.Find.Execute("85-43=42", $true, $false, $false, $false, $false, $true, 1, $false, "19+0=19", 2) | Out-Null
$d.Content.Find.Execute("12+52=64", $true, $false, $false, $false, $false, $true, 1, $false, "44-42=2", 2) | Out-Null
$d.Content.Find.Execute("58-3=55", $true, $false, $false, $false, $false, $true, 1, $false, "51-8=43", 2) | Out-Null
$d.Content.Find.Execute("68+30=98", $true, $false, $false, $false, $false, $true, 1, $false, "76-15=61", 2) | Out-Null
$d.Content.Find.Execute("41-20=21", $true, $false, $false, $false, $false, $true, 1, $false, "80-71=9", 2) | Out-Null
$d.Content.Find.Execute("74-12=62", $true, $false, $false, $false, $false, $true, 1, $false, "76-31=45", 2) | Out-Null
$d.Content.Find.Execute("1+87=88", $true, $false, $false, $false, $false, $true, 1, $false, "82-81=1", 2) | Out-Null
$d.Content.Find.Execute("28-3=25", $true, $false, $false, $false, $false, $true, 1, $false, "53+12=65", 2) | Out-Null
$d.Content.Find.Execute("27-11=16", $true, $false, $false, $false, $false, $true, 1, $false, "65-11=54", 2) | Out-Null
$d.Content.Find.Execute("72-7=65", $true, $false, $false, $false, $false, $true, 1, $false, "64+28=92", 2) | Out-Null
$d.Content.Find.Execute("87-12=75", $true, $false, $false, $false, $false, $true, 1, $false, "30-8=22", 2) | Out-Null
$d.Content.Find.Execute("67+24=91", $true, $false, $false, $false, $false, $true, 1, $false, "20-16=4", 2) | Out-Null
$d.Content.Find.Execute("20+71=91", $true, $false, $false, $false, $false, $true, 1, $false, "2+86=88", 2) | Out-Null
$d.Content.Find.Execute("88-11=77", $true, $false, $false, $false, $false, $true, 1, $false, "15-14=1", 2) | Out-Null
$d.Content.Find.Execute("85-33=52", $true, $false, $false, $false, $false, $true, 1, $false, "56+22=78", 2) | Out-Null
$d.Content.Find.Execute("74+21=95", $true, $false, $false, $false, $false, $true, 1, $false, "61-34=27", 2) | Out-Null
$d.Content.Find.Execute("4+52=56", $true, $false, $false, $false, $false, $true, 1, $false, "65+31=96", 2) | Out-Null
$d.Content.Find.Execute("4+95=99", $true, $false, $false, $false, $false, $true, 1, $false, "60+20=80", 2) | Out-Null
$d.Content.Find.Execute("72-32=40", $true, $false, $false, $false, $false, $true, 1, $false, "81-77=4", 2) | Out-Null
$d.Content.Find.Execute("58-0=58", $true, $false, $false, $false, $false, $true, 1, $false, "47-18=29", 2) | Out-Null
$d.Content.Find.Execute("16+10=26", $true, $false, $false, $false, $false, $true, 1, $false, "37+46=83", 2) | Out-Null
$d.Content.Find.Execute("40+55=95", $true, $false, $false, $false, $false, $true, 1, $false, "54-12=42", 2) | Out-Null
$d.Content.Find.Execute("28+45=73", $true, $false, $false, $false, $false, $true, 1, $false, "73-39=34", 2) | Out-Null
$d.Content.Find.Execute("88+8=96", $true, $false, $false, $false, $false, $true, 1, $false, "21+69=90", 2) | Out-Null
$d.Content.Find.Execute("70-1=69", $true, $false, $false, $false, $false, $true, 1, $false, "81-30=51", 2) | Out-Null
$d.Content.Find.Execute("34-24=10", $true, $false, $false, $false, $false, $true, 1, $false, "27+17=44", 2) | Out-Null
$d.Content.Find.Execute("21-9=12", $true, $false, $false, $false, $false, $true, 1, $false, "12+73=85", 2) | Out-Null
$d.Content.Find.Execute("44+4=48", $true, $false, $false, $false, $false, $true, 1, $false, "35-29=6", 2) | Out-Null
$d.Content.Find.Execute("89-69=20", $true, $false, $false, $false, $false, $true, 1, $false, "86+1=87", 2) | Out-Null
$d.Content.Find.Execute("95-16=79", $true, $false, $false, $false, $false, $true, 1, $false, "87-48=39", 2) | Out-Null
$d.Content.Find.Execute("25-17=8", $true, $false, $false, $false, $false, $true, 1, $false, "32-28=4", 2) | Out-Null
$d.Content.Find.Execute("71-47=24", $true, $false, $false, $false, $false, $true, 1, $false, "4+27=31", 2) | Out-Null
$d.Content.Find.Execute("59+30=89", $true, $false, $false, $false, $false, $true, 1, $false, "98-0=98", 2) | Out-Null
$d.Content.Find.Execute("15+33=48", $true, $false, $false, $false, $false, $true, 1, $false, "30+6=36", 2) | Out-Null
$d.Content.Find.Execute("49-15=34", $true, $false, $false, $false, $false, $true, 1, $false, "17+70=87", 2) | Out-Null
$d.Content.Find.Execute("21+24=45", $true, $false, $false, $false, $false, $true, 1, $false, "25-23=2", 2) | Out-Null
$d.Content.Find.Execute("97-13=84", $true, $false, $false, $false, $false, $true, 1, $false, "22+3=25", 2) | Out-Null
$d.Content.Find.Execute("39-0=39", $true, $false, $false, $false, $false, $true, 1, $false, "48+32=80", 2) | Out-Null
$d.Content.Find.Execute("50-6=44", $true, $false, $false, $false, $false, $true, 1, $false, "56+16=72", 2) | Out-Null
$d.Content.Find.Execute("70-32=38", $true, $false, $false, $false, $false, $true, 1, $false, "40-3=37", 2) | Out-Null
$d.Content.Find.Execute("96-21=75", $true, $false, $false, $false, $false, $true, 1, $false, "25+38=63", 2) | Out-Null
$d.Content.Find.Execute("96-3=93", $true, $false, $false, $false, $false, $true, 1, $false, "50+32=82", 2) | Out-Null
$d.Content.Find.Execute("99-82=17", $true, $false, $false, $false, $false, $true, 1, $false, "74-29=45", 2) | Out-Null
$d.Content.Find.Execute("12+50=62", $true, $false, $false, $false, $false, $true, 1, $false, "29+65=94", 2) | Out-Null
$d.Content.Find.Execute("9+24=33", $true, $false, $false, $false, $false, $true, 1, $false, "26+14=40", 2) | Out-Null
$d.Content.Find.Execute("83-38=45", $true, $false, $false, $false, $false, $true, 1, $false, "23+4=27", 2) | Out-Null
$d.Content.Find.Execute("68+6=74", $true, $false, $false, $false, $false, $true, 1, $false, "66+25=91", 2) | Out-Null
$d.Content.Find.Execute("80-45=35", $true, $false, $false, $false, $false, $true, 1, $false, "0+5=5", 2) | Out-Null
$d.Content.Find.Execute("25+68=93", $true, $false, $false, $false, $false, $true, 1, $false, "5+76=81", 2) | Out-Null
$d.Content.Find.Execute("25+37=62", $true, $false, $false, $false, $false, $true, 1, $false, "18+11=29", 2) | Out-Null
$d.Content.Find.Execute("93-10=83", $true, $false, $false, $false, $false, $true, 1, $false, "39-34=5", 2) | Out-Null
$d.Content.Find.Execute("93-13=80", $true, $false, $false, $false, $false, $true, 1, $false, "24+60=84", 2) | Out-Null
$d.Content.Find.Execute("38+18=56", $true, $false, $false, $false, $false, $true, 1, $false, "41-35=6", 2) | Out-Null
$d.Content.Find.Execute("34+33=67", $true, $false, $false, $false, $false, $true, 1, $false, "67-9=58", 2) | Out-Null
$d.Content.Find.Execute("45+29=74", $true, $false, $false, $false, $false, $true, 1, $false, "11+34=45", 2) | Out-Null
$d.Content.Find.Execute("1+12=13", $true, $false, $false, $false, $false, $true, 1, $false, "33+56=89", 2) | Out-Null
$d.Content.Find.Execute("15+55=70", $true, $false, $false, $false, $false, $true, 1, $false, "64+11=75", 2) | Out-Null
$d.Content.Find.Execute("42-25=17", $true, $false, $false, $false, $false, $true, 1, $false, "34+55=89", 2) | Out-Null
$d.Content.Find.Execute("71+16=87", $true, $false, $false, $false, $false, $true, 1, $false, "95-20=75", 2) | Out-Null
$d.Content.Find.Execute("44+23=67", $true, $false, $false, $false, $false, $true, 1, $false, "63-56=7", 2) | Out-Null
